$d = $word.ActiveDocument

# Locate the "Force Read" bullet paragraph (originally a single paragraph that
# the edit splits into two: the existing bullet, trimmed down, plus a new
# "Hint:" sub-bullet one indent level deeper).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*passing in a random num*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find the 'Force Read' paragraph to edit"
}

# Replace that whole paragraph with the new two-paragraph version. The leading
# runs (through the 'favoriteCharacters' endpoint mention) are reproduced
# verbatim from the source; only the trailing runs change, and a brand new
# "Hint" bullet (ListParagraph, ilvl 2, numId 2) follows it, carrying forward
# the relocated "_GoBack" bookmark and the "favoriteCharacters" spell-check
# tags exactly as Word would re-flag them.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="7965212F" w14:textId="099A6496" w:rsidR="00516F57" w:rsidRPr="00333A19" w:rsidRDefault="00516F57" w:rsidP="00886ECC"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r w:rsidRPr="00333A19"><w:t xml:space="preserve">The button “Force </w:t></w:r><w:r w:rsidR="00BB2176"><w:t>Read</w:t></w:r><w:r w:rsidRPr="00333A19"><w:t xml:space="preserve">” will GET the </w:t></w:r><w:r w:rsidR="0061161D" w:rsidRPr="00333A19"><w:t>data from the server endpoint</w:t></w:r><w:r w:rsidR="004B0A0A" w:rsidRPr="00333A19"><w:t xml:space="preserve"> /</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00333A19"><w:t>favoriteCharacters</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004B0A0A" w:rsidRPr="00333A19"><w:t>/{index</w:t></w:r><w:r w:rsidRPr="00333A19"><w:t xml:space="preserve">}. </w:t></w:r><w:r w:rsidRPr="00333A19"><w:t xml:space="preserve">The {index} will be a random index from the list on the server. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Hint: To get a random index, first make a normal GET request like you did for step c, then use the length of the resul</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>t to get a random index. Then immediately make a second request to /</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>favoriteCharacters</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/{index} and display those results and the index that was chosen.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$null = $targetPara.Range.InsertXML($xml)

Write-Output "Split the 'Force Read' bullet and added the new Hint sub-bullet."
